$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reshape the grid: drop the last 2 columns (I, J) and last 2 rows (11, 12) ---
$ws.Range("A10:I12").UnMerge()
$ws.Columns.Item(9).Delete()
$ws.Columns.Item(9).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(11).Delete()

# --- Rewrite all cell content for the new 6-alternative pairwise comparison matrix ---
$ws.Cells.Item(1, 2).Value = "Kia Rio"
$ws.Cells.Item(1, 3).Value = "Volkswagen Golf"
$ws.Cells.Item(1, 4).Value = "Toyota Corolla"
$ws.Cells.Item(1, 5).Value = "Skoda Octavia"
$ws.Cells.Item(1, 6).Value = "BMW 3 Series"
$ws.Cells.Item(1, 7).Value = "Hyundai Solaris"
$ws.Cells.Item(1, 8).Value = "Вектор приоритетов"
$ws.Cells.Item(2, 1).Value = "Kia Rio"
$ws.Cells.Item(2, 2).Value = "1"
$ws.Cells.Item(2, 3).Value = "3"
$ws.Cells.Item(2, 4).Value = "2"
$ws.Cells.Item(2, 5).Value = "2"
$ws.Cells.Item(2, 6).Value = "5"
$ws.Cells.Item(2, 7).Value = "1"
$ws.Cells.Item(2, 8).Value = "0.275"
$ws.Cells.Item(3, 1).Value = "Volkswagen Golf"
$ws.Cells.Item(3, 2).Value = "1/3"
$ws.Cells.Item(3, 3).Value = "1"
$ws.Cells.Item(3, 4).Value = "1/2"
$ws.Cells.Item(3, 5).Value = "1/2"
$ws.Cells.Item(3, 6).Value = "3"
$ws.Cells.Item(3, 7).Value = "1/3"
$ws.Cells.Item(3, 8).Value = "0.094"
$ws.Cells.Item(4, 1).Value = "Toyota Corolla"
$ws.Cells.Item(4, 2).Value = "1/2"
$ws.Cells.Item(4, 3).Value = "2"
$ws.Cells.Item(4, 4).Value = "1"
$ws.Cells.Item(4, 5).Value = "1"
$ws.Cells.Item(4, 6).Value = "4"
$ws.Cells.Item(4, 7).Value = "1/2"
$ws.Cells.Item(4, 8).Value = "0.156"
$ws.Cells.Item(5, 1).Value = "Skoda Octavia"
$ws.Cells.Item(5, 2).Value = "1/2"
$ws.Cells.Item(5, 3).Value = "2"
$ws.Cells.Item(5, 4).Value = "1"
$ws.Cells.Item(5, 5).Value = "1"
$ws.Cells.Item(5, 6).Value = "4"
$ws.Cells.Item(5, 7).Value = "1/2"
$ws.Cells.Item(5, 8).Value = "0.156"
$ws.Cells.Item(6, 1).Value = "BMW 3 Series"
$ws.Cells.Item(6, 2).Value = "1/5"
$ws.Cells.Item(6, 3).Value = "1/3"
$ws.Cells.Item(6, 4).Value = "1/4"
$ws.Cells.Item(6, 5).Value = "1/4"
$ws.Cells.Item(6, 6).Value = "1"
$ws.Cells.Item(6, 7).Value = "1/5"
$ws.Cells.Item(6, 8).Value = "0.044"
$ws.Cells.Item(7, 1).Value = "Hyundai Solaris"
$ws.Cells.Item(7, 2).Value = "1"
$ws.Cells.Item(7, 3).Value = "3"
$ws.Cells.Item(7, 4).Value = "2"
$ws.Cells.Item(7, 5).Value = "2"
$ws.Cells.Item(7, 6).Value = "5"
$ws.Cells.Item(7, 7).Value = "1"
$ws.Cells.Item(7, 8).Value = "0.275"
$ws.Cells.Item(8, 8).Value = "λ_max = 6.063"
$ws.Cells.Item(9, 8).Value = "ИС = 0.013"
$ws.Cells.Item(10, 8).Value = "ОС = 0.010"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 19.5
$ws.Columns.Item(2).ColumnWidth = 10.0
$ws.Columns.Item(3).ColumnWidth = 19.5
$ws.Columns.Item(4).ColumnWidth = 18.333333333333332
$ws.Columns.Item(5).ColumnWidth = 17.166666666666668
$ws.Columns.Item(6).ColumnWidth = 16.0
$ws.Columns.Item(7).ColumnWidth = 19.5
$ws.Columns.Item(8).ColumnWidth = 23.166666666666668

# --- Re-merge the summary block (lambda_max / CI / CR) ---
$ws.Range("A8:G10").Merge()

# --- Style: wrap text + vertical-center the data range ---
$rng = $ws.Range("A1:H10")
$rng.WrapText = $true
$rng.VerticalAlignment = -4108
